$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("household_size"), shifting the
# existing collector_name/payment_channel/.../national_id columns (and
# all their formatting) one column to the right. This mirrors adding a
# new "collector_id" field to the export.
$ws.Columns("C").Insert()

# The new column doesn't pick up an explicit width from its neighbours,
# so match it to column B (displayed width 17, which Excel's ColumnWidth
# property reports as 16.17 because of the internal padding offset).
$ws.Columns("C").ColumnWidth = 16.17

# Populate the new "collector_id" column: header + the two data rows.
$ws.Range("C1").Value = "collector_id"
$ws.Range("C2").Value = "IND-24-0000.0012"
$ws.Range("C3").Value = "IND-24-0000.0013"
